{"js": "// 1. Remove the old \"_GoBack\" bookmark (it sat right after the\n//    \"YEISON AGUIRRE OSORIO \" run on the cover page). Word always keeps\n//    only one \"_GoBack\" bookmark, marking the most recent edit location,\n//    so it has to be deleted before a new one is placed at the new edit\n//    spot further down in the document.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Locate the sentence that needs correcting, in the \"3.7. Paso 7:\"\n//    section: the maximum allowed number of iterations goes from 1000\n//    to 20000, and a clarifying sentence about the default (10000) is\n//    appended.\nconst body = context.document.body;\nconst searchResults = body.search(\n  \"entre 1 y 1000) que desea que realice el algoritmo.\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const target = searchResults.items[0];\n\n  // Replace the matched text with the same run formatting (Arial 12pt,\n  // es-CO), but split across four runs so the \"20\" insertion and the\n  // untouched \"000) que\" tail remain distinguishable, mirroring how the\n  // edit was actually typed.\n  const ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"es-CO\"/></w:rPr><w:t xml:space=\"preserve\">entre 1 y </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"es-CO\"/></w:rPr><w:t>20</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"es-CO\"/></w:rPr><w:t>000) que</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"es-CO\"/></w:rPr><w:t xml:space=\"preserve\"> desea que realice el algoritmo, si deja las iteraciones en 0 se realizaran 10000 iteraciones.</w:t></w:r></w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\n\n  target.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n\n  // 3. Word stamps a fresh \"_GoBack\" bookmark at the point of the last\n  //    edit on save - recreate it right after the sentence we just\n  //    fixed. Re-search for the trailing text (the original \"target\"\n  //    range does not grow to cover text inserted via insertOoxml) so\n  //    the bookmark anchors at the true end of the corrected sentence.\n  const tailResults = body.search(\"10000 iteraciones.\", { matchCase: true });\n  tailResults.load(\"items\");\n  await context.sync();\n\n  if (tailResults.items.length > 0) {\n    const tailEnd = tailResults.items[0].getRange(Word.RangeLocation.end);\n    tailEnd.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the old \"_GoBack\" bookmark (it sat right after the\n#    \"YEISON AGUIRRE OSORIO \" run on the cover page). Word keeps only\n#    one \"_GoBack\" bookmark - the most recent edit location - so the\n#    old one has to go before a new one appears at the new edit spot.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Locate the sentence that needs correcting in the \"3.7. Paso 7:\"\n#    section: the maximum allowed number of iterations goes from 1000\n#    to 20000, and a clarifying sentence about the default (10000) is\n#    appended.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"entre 1 y 1000) que desea que realice el algoritmo.\")\nif ($found) {\n    $rng.Text = \"entre 1 y 20000) que desea que realice el algoritmo, si deja las iteraciones en 0 se realizaran 10000 iteraciones.\"\n\n    # 3. Word stamps a fresh \"_GoBack\" bookmark at the point of the last\n    #    edit on save - recreate it right after the sentence we just\n    #    fixed.\n    $rng.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $rng)\n}\n"}
